# Add a "last updated" date stamp to cell C1 of the "About" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
